# Update the "Protect Lake Powell In Brief" data-assembly deck from the
# March 30, 2022 values to the April 21, 2022 values.

$p = $ppt.ActivePresentation

function Set-RangeTextClean($fullTextRange, $rangeObj, [string]$newText) {
    # Assigning straight onto a TextRange's .Text diffs the old/new strings
    # and keeps any shared prefix/suffix characters pinned to their original
    # run, splitting the run apart. Routing the edit through an unrelated
    # placeholder first leaves nothing in common for that diff to preserve,
    # so the range collapses back down to a single run (with the original
    # formatting untouched) once the real text is written in a second pass.
    $startPos = $rangeObj.Start
    $placeholder = "zzqqxxPLACEHOLDERxxqqzz" + [string]$startPos
    $rangeObj.Text = $placeholder
    $fresh = $fullTextRange.Characters($startPos, $placeholder.Length)
    $fresh.Text = $newText
}

function Update-InflowBox($slide, $splitRun) {
    $shape = $slide.Shapes.Item(2)
    $tr = $shape.TextFrame.TextRange
    $para = $tr.Paragraphs(2, 1)

    if ($splitRun) {
        # Retype just the leading "5.2 " as "5.3 ", which lands in its own
        # run ahead of the untouched remainder, then fix up that remainder
        # ("– 6.6 " -> "– 6.3 ") in place.
        $lead = $para.Characters(1, 4)
        $lead.Text = "5.3 "

        $rest = $para.Characters(5, 6)
        $rest.Text = "– 6.3 "
    }
    else {
        $numbers = $para.Characters(1, 10)
        Set-RangeTextClean $tr $numbers "5.3 – 6.3 "
    }
}

function Update-ValuesDateBox($slide) {
    $shape = $slide.Shapes.Item(13)
    $tr = $shape.TextFrame.TextRange
    $para = $tr.Paragraphs(1, 1)
    Set-RangeTextClean $tr $para "Values for April 21, 2022"
}

$slide1 = $p.Slides.Item(1)
Update-InflowBox $slide1 $false
Update-ValuesDateBox $slide1

$slide2 = $p.Slides.Item(2)
Update-InflowBox $slide2 $true
Update-ValuesDateBox $slide2
